$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new "Save" column in H1, matching the style used by the
# other header cells (e.g. G1). Copy G1's formatting to H1 first, then
# overwrite the value so the copied ("sum") text isn't kept.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Fill in the "Save" column values for rows 2-13.
$saveValues = @{
    2 = 0
    3 = 1
    4 = 0
    5 = 0
    6 = 1
    7 = 0
    8 = 1
    9 = 1
    10 = 0
    11 = 0
    12 = 1
    13 = 1
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
